$wb = $excel.ActiveWorkbook

# --- Sheet "Protocol (Server - Client)": mark the two open rows as Closed with a date ---
$wsProtocol = $wb.Worksheets.Item("Protocol (Server - Client)")
$wsProtocol.Range("C9").Value = "Closed"
$wsProtocol.Range("E9").Value = 41624
$wsProtocol.Range("C10").Value = "Closed"
$wsProtocol.Range("E10").Value = 41624

# --- Sheet "Simulator (Client)": tweak punctuation on three remarks ---
$wsSimulator = $wb.Worksheets.Item("Simulator (Client)")
$wsSimulator.Range("D20").Value = "Nog niet aangewerkt."
$wsSimulator.Range("D19").Value = "Ze rijden al wel alleen nog niet naar de juiste locatie."
$wsSimulator.Range("D10").Value = "Met een vector, die vanuit het protocol wordt gestuurd, werkt dit!"

# --- Reset each sheet's selection to B3 (leaves the scroll position at the top) ---
$wsController = $wb.Worksheets.Item("Controller (Server)")
[void]$wsController.Activate()
[void]$wsController.Range("B3").Select()

[void]$wsProtocol.Activate()
[void]$wsProtocol.Range("B3").Select()

[void]$wsSimulator.Activate()
[void]$wsSimulator.Range("B3").Select()

$wsManagement = $wb.Worksheets.Item("Management Interface (Mobiel)")
[void]$wsManagement.Activate()
[void]$wsManagement.Range("B3").Select()

# --- Make "Controller (Server)" the active tab shown when the workbook is reopened ---
[void]$wsController.Activate()
